$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BoM")

# Update footprint/value cells for the first component group (capacitors)
$ws.Range("F9").Value = "C_0402_1005Metric"

# Update the connector group's "Value" column entry
$ws.Range("E11").Value = "02x02"

# Update footprint for the resistor group
$ws.Range("F12").Value = "R_0402_1005Metric"

# Narrow the Footprint column (F) width (target stored width 40.7109375;
# Excel's ColumnWidth setter quantizes to whole-pixel boundaries, so we pick
# the input that lands on the closest achievable stored width)
$ws.Columns.Item(6).ColumnWidth = 39.8
